$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data edits ---

# B3: 2638 -> 3188
$ws.Range("B3").Value = 3188

# B9: 557 -> 1557
$ws.Range("B9").Value = 1557

# F23: 1105 -> 1205
$ws.Range("F23").Value = 1205

# Row 24 (Oct 22 entry): fill in bazar charge, note and meal counts
$ws.Range("F24").Value = 2375
$ws.Range("G24").Value = "225 taka baki"
$ws.Range("K24:T24").Value = 2

# --- New column (G) width ---
$ws.Columns.Item(7).ColumnWidth = 15.7

# --- View state ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 4
$win.ScrollRow = 13
$ws.Range("H34").Select()
